$d = $word.ActiveDocument

# Locate the unique "**" paragraph near the end of the document.
$findRng = $d.Content
$findRng.Find.Execute("**", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$targetStart = $findRng.Start

$count = $d.Paragraphs.Count
$starIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -eq $targetStart) {
        $starIndex = $i
        break
    }
}

$prevIndex = $starIndex - 1
$anchorPara = $d.Paragraphs.Item($prevIndex)
$starPara = $d.Paragraphs.Item($starIndex)

# Span covering the empty paragraph right before "**" through the "**" paragraph itself.
$combined = $d.Range($anchorPara.Range.Start, $starPara.Range.End)

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$xml = '<w:p ' + $ns + '><w:pPr><w:pStyle w:val="NormalnyWeb"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>KatDokumentyRodzaj</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>' + `
       '<w:p><w:pPr><w:pStyle w:val="NormalnyWeb"/></w:pPr></w:p>' + `
       '<w:p><w:pPr><w:pStyle w:val="NormalnyWeb"/></w:pPr></w:p>' + `
       '<w:p><w:pPr><w:pStyle w:val="NormalnyWeb"/></w:pPr></w:p>' + `
       '<w:p><w:pPr><w:pStyle w:val="NormalnyWeb"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>zMIANA</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> RODZAJU DOKUMENTU !!!!!!!</w:t></w:r></w:p>' + `
       '<w:p><w:pPr><w:pStyle w:val="NormalnyWeb"/></w:pPr><w:r><w:t>**</w:t></w:r></w:p>'

$combined.InsertXML($xml)
